# Apply the cryptos-list price/volume refresh described by the commit
# "Updated cryptos list on Thu May 25 16:46:23 UTC 2023 with GitHub Actions".
#
# All data cells on the sheet are plain text cells (inline strings), including
# the "Price" column which often holds digit-and-dot text such as "306.58" that
# Excel would otherwise auto-convert to a real number on assignment. Set-PriceText
# below works around that by assigning through a leading apostrophe (forces text)
# and then restoring the cell's original style so no stray formatting is introduced.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-PriceText($cell, [string]$value) {
    $origStyle = $cell.Style
    $cell.Value = "'" + $value
    $cell.Style = $origStyle
}

function Set-Text($cell, [string]$value) {
    $cell.Value = $value
}

Set-PriceText $ws.Cells.Item(2, 4) "26.440.33"
Set-Text $ws.Cells.Item(2, 5) "  +0.18%  "

Set-PriceText $ws.Cells.Item(3, 4) "1.811.15"
Set-Text $ws.Cells.Item(3, 5) "  +0.56%  "

Set-Text $ws.Cells.Item(4, 5) "  +0.18%  "

Set-Text $ws.Cells.Item(5, 5) "  -0.01%  "

Set-PriceText $ws.Cells.Item(6, 4) "306.58"
Set-Text $ws.Cells.Item(6, 5) "  -0.49%  "

Set-PriceText $ws.Cells.Item(7, 4) "0.4514"
Set-Text $ws.Cells.Item(7, 5) "  +0.01%  "

Set-PriceText $ws.Cells.Item(8, 4) "0.3589"
Set-Text $ws.Cells.Item(8, 5) "  -1.48%  "

Set-PriceText $ws.Cells.Item(9, 4) "46.51"
Set-Text $ws.Cells.Item(9, 5) "  +3.94%  "

Set-Text $ws.Cells.Item(10, 5) "  -0.18%  "

Set-PriceText $ws.Cells.Item(11, 4) "0.8929"
Set-Text $ws.Cells.Item(11, 5) "  +3.01%  "

Set-PriceText $ws.Cells.Item(12, 4) "0.07821"
Set-Text $ws.Cells.Item(12, 5) "  +0.59%  "

Set-PriceText $ws.Cells.Item(13, 4) "19.39"
Set-Text $ws.Cells.Item(13, 5) "  +0.88%  "

Set-PriceText $ws.Cells.Item(14, 4) "1.816.36"
Set-Text $ws.Cells.Item(14, 5) "  +0.27%  "

Set-PriceText $ws.Cells.Item(15, 4) "5.290"

Set-PriceText $ws.Cells.Item(16, 4) "6.317"
Set-Text $ws.Cells.Item(16, 5) "  +0.16%  "

Set-PriceText $ws.Cells.Item(17, 4) "85.14"
Set-Text $ws.Cells.Item(17, 5) "  -1.07%  "

Set-PriceText $ws.Cells.Item(19, 4) "0.000008517"
Set-Text $ws.Cells.Item(19, 5) "  -0.25%  "

Set-PriceText $ws.Cells.Item(20, 4) "1.008"
Set-Text $ws.Cells.Item(20, 5) "  +0.01%  "

Set-PriceText $ws.Cells.Item(21, 4) "26.484.04"
Set-Text $ws.Cells.Item(21, 5) "  +0.21%  "

Set-Text $ws.Cells.Item(22, 5) "  +0.06%  "

Set-PriceText $ws.Cells.Item(23, 4) "4.972"
Set-Text $ws.Cells.Item(23, 5) "  +0.34%  "

Set-PriceText $ws.Cells.Item(24, 4) "2.030.43"
Set-Text $ws.Cells.Item(24, 5) "  +0.11%  "

Set-PriceText $ws.Cells.Item(25, 4) "10.52"
Set-Text $ws.Cells.Item(25, 5) "  +1.06%  "

Set-PriceText $ws.Cells.Item(26, 4) "1.959"
Set-Text $ws.Cells.Item(26, 5) "  -0.78%  "

Set-PriceText $ws.Cells.Item(27, 4) "151.95"
Set-Text $ws.Cells.Item(27, 5) "  +1.52%  "

Set-Text $ws.Cells.Item(28, 5) "  -0.25%  "

Set-PriceText $ws.Cells.Item(29, 4) "2.064"
Set-Text $ws.Cells.Item(29, 5) "  +4.48%  "

Set-PriceText $ws.Cells.Item(30, 4) "112.26"
Set-Text $ws.Cells.Item(30, 5) "  -0.07%  "

Set-PriceText $ws.Cells.Item(31, 4) "4.859"
Set-Text $ws.Cells.Item(31, 5) "  +0.29%  "

Set-PriceText $ws.Cells.Item(32, 4) "0.08694"
Set-Text $ws.Cells.Item(32, 5) "  +0.79%  "

Set-PriceText $ws.Cells.Item(33, 4) "3.118"
Set-Text $ws.Cells.Item(33, 5) "  +2.96%  "

Set-PriceText $ws.Cells.Item(34, 4) "2.793"
Set-Text $ws.Cells.Item(34, 5) "  +9.90%  "

Set-Text $ws.Cells.Item(35, 2) "ImmutableX"
Set-Text $ws.Cells.Item(35, 3) "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-PriceText $ws.Cells.Item(35, 4) "0.7315"
Set-Text $ws.Cells.Item(35, 5) "  +0.75%  "

Set-Text $ws.Cells.Item(36, 2) "Filecoin"
Set-Text $ws.Cells.Item(36, 3) "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-PriceText $ws.Cells.Item(36, 4) "4.458"
Set-Text $ws.Cells.Item(36, 5) "  +0.78%  "

Set-Text $ws.Cells.Item(37, 5) "  -0.06%  "

Set-PriceText $ws.Cells.Item(38, 4) "1.077"
Set-Text $ws.Cells.Item(38, 5) "  +0.66%  "

Set-PriceText $ws.Cells.Item(39, 4) "0.01926"
Set-Text $ws.Cells.Item(39, 5) "  +0.76%  "

Set-PriceText $ws.Cells.Item(40, 4) "0.05121"
Set-Text $ws.Cells.Item(40, 5) "  +1.35%  "

Set-PriceText $ws.Cells.Item(41, 4) "2.898"
Set-Text $ws.Cells.Item(41, 5) "  +0.93%  "

Set-PriceText $ws.Cells.Item(42, 4) "0.5092"
Set-Text $ws.Cells.Item(42, 5) "  +4.10%  "

Set-PriceText $ws.Cells.Item(43, 4) "6.780"
Set-Text $ws.Cells.Item(43, 5) "  -2.76%  "

Set-Text $ws.Cells.Item(44, 5) "  -3.14%  "

Set-PriceText $ws.Cells.Item(45, 4) "8.054"
Set-Text $ws.Cells.Item(45, 5) "  -0.52%  "

Set-Text $ws.Cells.Item(46, 5) "  -0.07%  "

Set-PriceText $ws.Cells.Item(47, 4) "0.4666"
Set-Text $ws.Cells.Item(47, 5) "  +1.46%  "

Set-PriceText $ws.Cells.Item(48, 4) "10.01"
Set-Text $ws.Cells.Item(48, 5) "  +0.72%  "

Set-PriceText $ws.Cells.Item(49, 4) "100.00"
Set-Text $ws.Cells.Item(49, 5) "  -1.39%  "

Set-PriceText $ws.Cells.Item(50, 4) "1.573"
Set-Text $ws.Cells.Item(50, 5) "  -0.31%  "

Set-PriceText $ws.Cells.Item(51, 4) "0.05990"
Set-Text $ws.Cells.Item(51, 5) "  +0.03%  "
